$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 386450.56
$ws.Range("J17").Value = 386450.56
$ws.Range("L17").Value = 1159351.68
$ws.Range("N17").Value = -1159687.68
$ws.Range("H100").Value = 3856.2856
$ws.Range("I100").Value = 2676.5557
$ws.Range("K100").Value = 2676.5557
$ws.Range("M100").Value = -2135.5557
$ws.Range("H111").Value = 1081.3636
$ws.Range("I111").Value = 1130.8572
$ws.Range("K111").Value = 3392.5716
$ws.Range("M111").Value = -325.5715999999998
$ws.Range("H112").Value = 1110.421
$ws.Range("J112").Value = 1118.2354
$ws.Range("L112").Value = 3354.7062
$ws.Range("N112").Value = -5570.706200000001
$ws.Range("H115").Value = 3417.25
$ws.Range("I115").Value = 3417.25
$ws.Range("K115").Value = 10251.75
$ws.Range("M115").Value = -8684.75
$ws.Range("H125").Value = 1935
$ws.Range("I125").Value = 1537.3334
$ws.Range("J125").Value = 2332.6667
$ws.Range("K125").Value = 13836.0006
$ws.Range("L125").Value = 20994.0003
$ws.Range("M125").Value = -11376.0006
$ws.Range("N125").Value = -25914.0003
$ws.Range("H129").Value = 1841.5
$ws.Range("I129").Value = 679
$ws.Range("J129").Value = 4399
$ws.Range("K129").Value = 2037
$ws.Range("L129").Value = 13197
$ws.Range("M129").Value = 2963
$ws.Range("N129").Value = -23197
$ws.Range("H132").Value = 2399.8057
$ws.Range("I132").Value = 2254.4546
$ws.Range("J132").Value = 3998.6667
$ws.Range("K132").Value = 6763.3638
$ws.Range("L132").Value = 11996.0001
$ws.Range("M132").Value = -4233.3638
$ws.Range("N132").Value = -17056.0001
$ws.Range("H137").Value = 2015.4
$ws.Range("J137").Value = 2583.6667
$ws.Range("L137").Value = 7751.000100000001
$ws.Range("N137").Value = -12851.0001
$ws.Range("H138").Value = 2544.8462
$ws.Range("J138").Value = 2211.4075
$ws.Range("L138").Value = 6634.2225
$ws.Range("N138").Value = -16914.2225

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25602.611
$ws.Range("I32").Value = 4849.8604
$ws.Range("J32").Value = 174330.67
$ws.Range("K32").Value = 4849.8604
$ws.Range("L32").Value = 174330.67
$ws.Range("M32").Value = -4562.8604
$ws.Range("N32").Value = -174904.67
$ws.Range("H74").Value = 1613.4791
$ws.Range("I74").Value = 1380.2142
$ws.Range("J74").Value = 3246.3333
$ws.Range("K74").Value = 1380.2142
$ws.Range("L74").Value = 3246.3333
$ws.Range("M74").Value = -506.2141999999999
$ws.Range("N74").Value = -4994.3333
$ws.Range("H77").Value = 1613.4791
$ws.Range("I77").Value = 1380.2142
$ws.Range("J77").Value = 3246.3333
$ws.Range("K77").Value = 6901.071
$ws.Range("L77").Value = 16231.6665
$ws.Range("M77").Value = -2533.071
$ws.Range("N77").Value = -24967.6665
$ws.Range("H122").Value = 2952.9
$ws.Range("I122").Value = 2957.3333
$ws.Range("K122").Value = 8871.999899999999
$ws.Range("M122").Value = -6421.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 680
$ws.Range("I22").Value = 818.75
$ws.Range("J22").Value = 125
$ws.Range("K22").Value = 818.75
$ws.Range("L22").Value = 125
$ws.Range("M22").Value = -645.75
$ws.Range("N22").Value = -471
$ws.Range("H86").Value = 1970.8572
$ws.Range("I86").Value = 1666
$ws.Range("J86").Value = 2580.5715
$ws.Range("K86").Value = 1666
$ws.Range("L86").Value = 2580.5715
$ws.Range("M86").Value = -543
$ws.Range("N86").Value = -4826.5715
$ws.Range("H89").Value = 1970.8572
$ws.Range("I89").Value = 1666
$ws.Range("J89").Value = 2580.5715
$ws.Range("K89").Value = 8330
$ws.Range("L89").Value = 12902.8575
$ws.Range("M89").Value = -2714
$ws.Range("N89").Value = -24134.8575
$ws.Range("H94").Value = 4373.476
$ws.Range("I94").Value = 4102.4443
$ws.Range("K94").Value = 4102.4443
$ws.Range("M94").Value = -3651.4443
$ws.Range("H107").Value = 93581
$ws.Range("I107").Value = 127175
$ws.Range("J107").Value = 3997
$ws.Range("K107").Value = 127175
$ws.Range("L107").Value = 3997
$ws.Range("M107").Value = -125255
$ws.Range("N107").Value = -7837
$ws.Range("H132").Value = 129011.5
$ws.Range("J132").Value = 129011.5
$ws.Range("L132").Value = 129011.5
$ws.Range("N132").Value = -139131.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 13282.167
$ws.Range("I62").Value = 16652.777
$ws.Range("J62").Value = 3170.3333
$ws.Range("K62").Value = 16652.777
$ws.Range("L62").Value = 3170.3333
$ws.Range("M62").Value = -16028.777
$ws.Range("N62").Value = -4418.3333
$ws.Range("H65").Value = 13282.167
$ws.Range("I65").Value = 16652.777
$ws.Range("J65").Value = 3170.3333
$ws.Range("K65").Value = 83263.88499999999
$ws.Range("L65").Value = 15851.6665
$ws.Range("M65").Value = -80143.88499999999
$ws.Range("N65").Value = -22091.6665
$ws.Range("H132").Value = 2598.1
$ws.Range("I132").Value = 2710.2188
$ws.Range("J132").Value = 2398.7778
$ws.Range("K132").Value = 8130.6564
$ws.Range("L132").Value = 7196.3334
$ws.Range("M132").Value = -5600.6564
$ws.Range("N132").Value = -12256.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3258
$ws.Range("I64").Value = 1459.75
$ws.Range("J64").Value = 4285.5713
$ws.Range("K64").Value = 4379.25
$ws.Range("L64").Value = 12856.7139
$ws.Range("M64").Value = -4109.25
$ws.Range("N64").Value = -13396.7139
$ws.Range("H67").Value = 3258
$ws.Range("I67").Value = 1459.75
$ws.Range("J67").Value = 4285.5713
$ws.Range("K67").Value = 4379.25
$ws.Range("L67").Value = 12856.7139
$ws.Range("M67").Value = -3443.25
$ws.Range("N67").Value = -14728.7139
$ws.Range("H98").Value = 1383.1666
$ws.Range("J98").Value = 1449
$ws.Range("L98").Value = 4347
$ws.Range("N98").Value = -7343
$ws.Range("H136").Value = 1007
$ws.Range("I136").Value = 1007
$ws.Range("K136").Value = 3021
$ws.Range("M136").Value = 2079
$ws.Range("H137").Value = 3124.625
$ws.Range("J137").Value = 4300.857
$ws.Range("L137").Value = 12902.571
$ws.Range("N137").Value = -23102.571
$ws.Range("H141").Value = 2808
$ws.Range("I141").Value = 2543.2222
$ws.Range("J141").Value = 3999.5
$ws.Range("K141").Value = 7629.6666
$ws.Range("L141").Value = 11998.5
$ws.Range("M141").Value = -2449.6666
$ws.Range("N141").Value = -22358.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3355.8333
$ws.Range("I80").Value = 2713.5
$ws.Range("J80").Value = 4640.5
$ws.Range("K80").Value = 2713.5
$ws.Range("L80").Value = 4640.5
$ws.Range("M80").Value = -1715.5
$ws.Range("N80").Value = -6636.5
$ws.Range("H83").Value = 3355.8333
$ws.Range("I83").Value = 2713.5
$ws.Range("J83").Value = 4640.5
$ws.Range("K83").Value = 13567.5
$ws.Range("L83").Value = 23202.5
$ws.Range("M83").Value = -8575.5
$ws.Range("N83").Value = -33186.5
$ws.Range("H122").Value = 1507.5385
$ws.Range("I122").Value = 1419.4445
$ws.Range("J122").Value = 1705.75
$ws.Range("K122").Value = 4258.333500000001
$ws.Range("L122").Value = 5117.25
$ws.Range("M122").Value = -1808.333500000001
$ws.Range("N122").Value = -10017.25
$ws.Range("H132").Value = 9666.223
$ws.Range("I132").Value = 8499.5
$ws.Range("J132").Value = 11999.667
$ws.Range("K132").Value = 25498.5
$ws.Range("L132").Value = 35999.001
$ws.Range("M132").Value = -22968.5
$ws.Range("N132").Value = -41059.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6609.875
$ws.Range("J7").Value = 5499.875
$ws.Range("L7").Value = 5499.875
$ws.Range("N7").Value = -5723.875
$ws.Range("H40").Value = 5499.8423
$ws.Range("J40").Value = 9332.6
$ws.Range("L40").Value = 9332.6
$ws.Range("N40").Value = -9604.6
$ws.Range("H46").Value = 28891.5
$ws.Range("J46").Value = 3204.1
$ws.Range("L46").Value = 3204.1
$ws.Range("N46").Value = -3580.1
$ws.Range("H61").Value = 108487.82
$ws.Range("I61").Value = 114377
$ws.Range("K61").Value = 114377
$ws.Range("M61").Value = -114175
$ws.Range("H68").Value = 1796
$ws.Range("I68").Value = 1796
$ws.Range("K68").Value = 1796
$ws.Range("M68").Value = -1047
$ws.Range("H71").Value = 1796
$ws.Range("I71").Value = 1796
$ws.Range("K71").Value = 8980
$ws.Range("M71").Value = -5236
$ws.Range("H100").Value = 17859.92
$ws.Range("I100").Value = 2844.2222
$ws.Range("K100").Value = 2844.2222
$ws.Range("M100").Value = -2303.2222
$ws.Range("H113").Value = 108487.82
$ws.Range("I113").Value = 114377
$ws.Range("K113").Value = 114377
$ws.Range("M113").Value = -112207
$ws.Range("H122").Value = 7807.56
$ws.Range("I122").Value = 9230.6875
$ws.Range("K122").Value = 27692.0625
$ws.Range("M122").Value = -25242.0625
$ws.Range("H126").Value = 6609.875
$ws.Range("J126").Value = 5499.875
$ws.Range("L126").Value = 16499.625
$ws.Range("N126").Value = -21439.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2479.7778
$ws.Range("I100").Value = 2479.7778
$ws.Range("K100").Value = 4959.5556
$ws.Range("M100").Value = -4418.5556
$ws.Range("H113").Value = 384.875
$ws.Range("I113").Value = 206.8
$ws.Range("K113").Value = 620.4000000000001
$ws.Range("M113").Value = 1549.6
$ws.Range("H132").Value = 4004.0889
$ws.Range("I132").Value = 4292
$ws.Range("K132").Value = 12876
$ws.Range("M132").Value = -10346
